# Swap the "step" contents of the TC3 and TC4 test case blocks.
# TC3 (row 21 label) currently holds the "detalhar diária" step in row 25.
# TC4 (row 28 label) currently holds the "analisar prestação de contas" step in row 32.
# After the edit, TC3's step (row 25) should be "analisar prestação de contas"
# and TC4's step (row 32) should be "detalhar diária" - the TC3/TC4 labels
# themselves stay where they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldB25 = $ws.Range("B25").Value2
$oldD25 = $ws.Range("D25").Value2
$oldB32 = $ws.Range("B32").Value2
$oldD32 = $ws.Range("D32").Value2

$ws.Range("B25").Value2 = $oldB32
$ws.Range("D25").Value2 = $oldD32
$ws.Range("B32").Value2 = $oldB25
$ws.Range("D32").Value2 = $oldD25
